# Insert a new data row at row 189 (pushing existing rows 189:287 down to 190:288)
# then populate the new row with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(189).Insert()

$ws.Cells.Item(189, 1).Value = 3
$ws.Cells.Item(189, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(189, 3).Value = "Coquimbo"
$ws.Cells.Item(189, 4).Value = 44572
$ws.Cells.Item(189, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(189, 5).Value = 5
$ws.Cells.Item(189, 6).Value = 100112031
$ws.Cells.Item(189, 7).Value = "Poroto verde"
$ws.Cells.Item(189, 8).Value = "Magnum"
$ws.Cells.Item(189, 9).Value = "Primera"
$ws.Cells.Item(189, 10).Value = 73
$ws.Cells.Item(189, 11).Value = 27000
$ws.Cells.Item(189, 12).Value = 28000
$ws.Cells.Item(189, 13).Value = 27521
$ws.Cells.Item(189, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(189, 15).Value = "Provincia de Talca"
$ws.Cells.Item(189, 16).Value = 1101
$ws.Cells.Item(189, 17).Value = 25
$ws.Cells.Item(189, 18).Value = "Hortaliza"
